# The workbook's shared-strings table contained a stray "." placeholder
# string. Re-entering the correct values for the cells that pointed at it
# removes the placeholder from xl/sharedStrings.xml (uniqueCount 57 -> 56)
# and every other shared-string index shifts down by one accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 ("Chili Ginger Tofu") serving size was mis-recorded as "." -> "1 cup"
$ws.Range("B16").Value = "1 cup"

# The "Vegan" column ("N") cells that were mis-recorded as "."
$ws.Range("G7").Value = "N"
$ws.Range("G13").Value = "N"
$ws.Range("G15").Value = "N"

# Restore the last-saved selection/view state
$ws.Range("E42").Select()
